$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial value (45204 -> 2023-10-05).
# Update every data row (2 through 134) to the new value 45205 (2023-10-06).
$ws.Range("C2:C134").Value = 45205
